$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing "NetNratio" header to "Net_N2Ar_ratio"
$ws.Range("N1").Value = "Net_N2Ar_ratio"

# New headers for the N-saturation-ratio columns
$ws.Range("O1").Value = "uM N2"
$ws.Range("P1").Value = "uM N2 field expected"
$ws.Range("Q1").Value = "NsatRatio"

# Fill O (uM N2 = N2/Ar ratio * uM Ar exp equilibrium), P (uM N2 field expected,
# copied from the uM N2 exp-equilibrium column) and Q (= O / P) for every data row.
for ($r = 2; $r -le 19; $r++) {
    $nToArRatio = $ws.Cells.Item($r, 4).Value2   # column D: N2/Ar
    $arExpEq    = $ws.Cells.Item($r, 6).Value2   # column F: uM Ar (exp equilibrium)
    $n2ExpEq    = $ws.Cells.Item($r, 7).Value2   # column G: uM N2 (exp equilibrium)

    $ws.Cells.Item($r, 15).Value = $nToArRatio * $arExpEq   # O: uM N2
    $ws.Cells.Item($r, 16).Value = $n2ExpEq                 # P: uM N2 field expected
    $ws.Cells.Item($r, 17).Formula = "=O$r/P$r"             # Q: NsatRatio
}

# Best-effort column width updates (engine quantizes to 1/6-character steps,
# so these land as close as the interop layer allows to the authored widths).
$ws.Columns.Item(14).ColumnWidth = 13.6666666666667
$ws.Range("P:Q").ColumnWidth = 18.1666666666667

# Match the author's final selection
$ws.Range("H32").Select()
